# Add a "date" column (C) next to the existing name/age columns, mirroring
# the "Add dates to XLS and SPSS tests" fixture change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column.
$ws.Range("C1").Value = "date"

# Give the new date cells an explicit yyyy-mm-dd display format *before*
# writing the values, so the engine doesn't synthesize its own default
# date format (e.g. m/d/yy) from the auto-detected date literals.
$ws.Range("C2:C5").NumberFormat = "yyyy\-mm\-dd"

$ws.Range("C2").Value = "1985-01-01"
$ws.Range("C3").Value = "1990-01-01"
$ws.Range("C4").Value = "2010-07-07"
$ws.Range("C5").Value = "1993-01-01"

# Keep the selection where the sheet was left after the edit.
[void]$ws.Range("J18").Select()
